$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.325.60"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "'2.522.78"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'316.32"
$ws.Range("E5").Value = "  +3.41%  "
$ws.Range("D6").Value = "'93.59"
$ws.Range("E6").Value = "  -6.38%  "
$ws.Range("D7").Value = "'0.571"
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.525"
$ws.Range("E9").Value = "  -4.05%  "
$ws.Range("D10").Value = "'35.34"
$ws.Range("E10").Value = "  -5.58%  "
$ws.Range("D11").Value = "'0.0803"
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("D12").Value = "'7.56"
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").Value = "'2.908.53"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'15.28"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'2.463.62"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").Value = "'0.838"
$ws.Range("D18").Value = "'42.392.86"
$ws.Range("D19").Value = "'12.83"
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("D20").Value = "'6.53"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").Value = "'0.0₃0954"
$ws.Range("E21").Value = "  -3.93%  "
$ws.Range("D22").Value = "'70.47"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").Value = "'249.29"
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("D24").Value = "'2.94"
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("D25").Value = "'2.00"
$ws.Range("E25").Value = "  -3.39%  "
$ws.Range("D26").Value = "'26.18"
$ws.Range("E26").Value = "  -4.91%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "'2.35"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("D29").Value = "'10.09"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").Value = "'38.72"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("D31").Value = "'5.88"
$ws.Range("E31").Value = "  -5.80%  "
$ws.Range("D32").Value = "'156.18"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("D33").Value = "'19.36"
$ws.Range("E33").Value = "  +4.45%  "
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("D35").Value = "'3.27"
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("D37").Value = "'0.0778"
$ws.Range("E37").Value = "  -2.78%  "
$ws.Range("D38").Value = "'0.110"
$ws.Range("E38").Value = "  -4.75%  "
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("D40").Value = "'23.58"
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("D41").Value = "'2.31"
$ws.Range("E41").Value = "  +10.56%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").Value = "'3.76"
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0298"
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'3.27"
$ws.Range("E45").Value = "  -6.22%  "
$ws.Range("D46").Value = "'2.007.14"
$ws.Range("E46").Value = "  -2.98%  "
$ws.Range("D47").Value = "'84.06"
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("D48").Value = "'8.78"
$ws.Range("E48").Value = "  -3.20%  "
$ws.Range("D49").Value = "'2.762.02"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("D50").Value = "'72.60"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").Value = "'101.43"
$ws.Range("E51").Value = "  -2.31%  "
